$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows 4:41 down to 5:42.
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new data record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44670
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112040
$ws.Range("G4").Value = "Cilantro"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 550
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 575
$ws.Range("N4").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 575
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
